$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.046.39'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '1.646.80'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.82%  '
$ws.Range("D5").Value = "'216.73"
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = "'0.507"
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("D9").Value = "'0.0641"
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("D10").Value = "'19.68"
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '1.875.90'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.30"
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("D14").Value = '1.624.84'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").Value = '26.048.52'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").Value = "'193.36"
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = "'4.36"
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").Value = "'9.96"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = "'6.25"
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'1.82"
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").Value = "'0.132"
$ws.Range("E25").Value = '  +5.68%  '
$ws.Range("D26").Value = "'144.46"
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("D28").Value = "'6.94"
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("E30").Value = '  +1.46%  '
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").Value = "'3.29"
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("D36").Value = "'0.906"
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = '1.131.32'
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").Value = "'0.542"
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").Value = "'2.48"
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("E40").Value = '  +0.94%  '
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("D42").Value = "'99.54"
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '1.784.68'
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("E45").Value = '  +4.16%  '
$ws.Range("D46").Value = "'56.77"
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("D47").Value = "'0.0529"
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("E49").Value = '  +1.86%  '
$ws.Range("D50").Value = "'0.416"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = '  -0.22%  '
